$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 4) down to the new row 5
$ws.Range("A4:H4").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122)

$ws.Cells.Item(5, 1).Value = "2025-08-12 13:07:51 UTC"
$ws.Cells.Item(5, 2).Value = "2025-08-12 18:37:51 IST"
$ws.Cells.Item(5, 3).Value = "SKIPPED"
$ws.Cells.Item(5, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item(5, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$ws.Cells.Item(5, 6).Value = ""
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = ""
